$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61, shifting existing rows 61-69 down to 62-70
$ws.Rows.Item(61).Insert()

# Populate the new row 61 with the inserted record's data
$ws.Cells.Item(61, 1).Value = 5
$ws.Cells.Item(61, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(61, 3).Value = "Maule"
$ws.Cells.Item(61, 4).Value = 44841
$ws.Cells.Item(61, 4).NumberFormat = $ws.Cells.Item(62, 4).NumberFormat
$ws.Cells.Item(61, 5).Value = 7
$ws.Cells.Item(61, 6).Value = 300000000
$ws.Cells.Item(61, 7).Value = "Espárragos"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 2500
$ws.Cells.Item(61, 11).Value = 1000
$ws.Cells.Item(61, 12).Value = 1000
$ws.Cells.Item(61, 13).Value = 1000
$ws.Cells.Item(61, 14).Value = "`$/kilo"
$ws.Cells.Item(61, 15).Value = "Provincia de Linares"
$ws.Cells.Item(61, 16).Value = 1000
$ws.Cells.Item(61, 17).Value = 1
$ws.Cells.Item(61, 18).Value = "Hortaliza"
